# Updated symbol list on Mon Jan 23 21:24:38 UTC 2023 with GitHub Actions
# Refresh Price (column D) and Volume(1h) (column E) figures for the
# coinranking.com crypto snapshot on the active worksheet.
#
# The cells hold plain text (e.g. "306.05", "1.39%") rather than real
# numbers/percentages, so each value is forced to Text format before being
# written and the style is reset back to Normal afterwards so the cell
# keeps its original (unstyled) appearance.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRange, $value) {
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $value
    $cellRange.Style = "Normal"
}

Set-TextValue $ws.Range("D2")  "306.05"
Set-TextValue $ws.Range("E2")  "1.39%"

Set-TextValue $ws.Range("D3")  "36.20"
Set-TextValue $ws.Range("E3")  "-1.84%"

Set-TextValue $ws.Range("D4")  "5.050"
Set-TextValue $ws.Range("E4")  "0.88%"

Set-TextValue $ws.Range("D5")  "0.07927"
Set-TextValue $ws.Range("E5")  "3.34%"

Set-TextValue $ws.Range("D6")  "2.240"
Set-TextValue $ws.Range("E6")  "8.83%"

Set-TextValue $ws.Range("D7")  "8.006"
Set-TextValue $ws.Range("E7")  "0.46%"

Set-TextValue $ws.Range("D8")  "0.9283"
Set-TextValue $ws.Range("E8")  "1.30%"

Set-TextValue $ws.Range("D9")  "0.09824"
Set-TextValue $ws.Range("E9")  "3.56%"

Set-TextValue $ws.Range("D10") "0.1869"
Set-TextValue $ws.Range("E10") "0.87%"

Set-TextValue $ws.Range("D11") "0.09003"

Set-TextValue $ws.Range("D12") "0.03719"
Set-TextValue $ws.Range("E12") "3.24%"

Set-TextValue $ws.Range("D13") "0.09936"
Set-TextValue $ws.Range("E13") "-0.62%"

Set-TextValue $ws.Range("D14") "0.001442"
Set-TextValue $ws.Range("E14") "-2.25%"

Set-TextValue $ws.Range("D15") "0.005641"
Set-TextValue $ws.Range("E15") "-1.69%"

Set-TextValue $ws.Range("D16") "3.467"
Set-TextValue $ws.Range("E16") "-0.14%"

Set-TextValue $ws.Range("D17") "4.151"
Set-TextValue $ws.Range("E17") "3.10%"

Set-TextValue $ws.Range("D18") "2.632"
Set-TextValue $ws.Range("E18") "10.63%"

Set-TextValue $ws.Range("E19") "0.05%"

Set-TextValue $ws.Range("E20") "-1.03%"

Set-TextValue $ws.Range("D21") "5.069"
Set-TextValue $ws.Range("E21") "2.27%"

Set-TextValue $ws.Range("E22") "1.49%"

Set-TextValue $ws.Range("D23") "0.04559"
Set-TextValue $ws.Range("E23") "-1.10%"

Set-TextValue $ws.Range("E24") "-0.21%"

Set-TextValue $ws.Range("D25") "0.004783"
Set-TextValue $ws.Range("E25") "-5.93%"

Set-TextValue $ws.Range("D26") "0.0001299"
Set-TextValue $ws.Range("E26") "-7.97%"

Set-TextValue $ws.Range("D39") "0.01907"
Set-TextValue $ws.Range("E39") "9.28%"

Set-TextValue $ws.Range("D40") "0.04907"
Set-TextValue $ws.Range("E40") "6.89%"

Set-TextValue $ws.Range("D41") "0.007827"
Set-TextValue $ws.Range("E41") "1.50%"

Set-TextValue $ws.Range("D42") "0.1395"
Set-TextValue $ws.Range("E42") "0.40%"

Set-TextValue $ws.Range("D43") "0.007807"
Set-TextValue $ws.Range("E43") "0.41%"

Set-TextValue $ws.Range("D44") "0.002124"
Set-TextValue $ws.Range("E44") "-2.22%"

Set-TextValue $ws.Range("D45") "0.01143"
Set-TextValue $ws.Range("E45") "10.36%"

Set-TextValue $ws.Range("D46") "0.00006287"
Set-TextValue $ws.Range("E46") "-0.34%"

Set-TextValue $ws.Range("D47") "0.00000000750"
Set-TextValue $ws.Range("E47") "-0.85%"

Set-TextValue $ws.Range("E48") "49.46%"

Set-TextValue $ws.Range("D49") "0.001801"
Set-TextValue $ws.Range("E49") "-10.69%"

Set-TextValue $ws.Range("D50") "0.00002099"
Set-TextValue $ws.Range("E50") "-0.85%"

Set-TextValue $ws.Range("D51") "0.0001999"
Set-TextValue $ws.Range("E51") "-0.85%"
